$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1429.4445
$ws.Range("J17").Value = 1429.4445
$ws.Range("L17").Value = 4288.333500000001
$ws.Range("N17").Value = -4624.333500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 58824652
$ws.Range("I42").Value = 76924060
$ws.Range("J42").Value = 1568.75
$ws.Range("K42").Value = 230772180
$ws.Range("L42").Value = 4706.25
$ws.Range("M42").Value = -230771950
$ws.Range("N42").Value = -5166.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3833.3333
$ws.Range("I51").Value = 3187.5
$ws.Range("K51").Value = 3187.5
$ws.Range("M51").Value = -2703.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1397
$ws.Range("I132").Value = 1080.3334
$ws.Range("J132").Value = 2663.6667
$ws.Range("K132").Value = 3241.0002
$ws.Range("L132").Value = 7991.000100000001
$ws.Range("M132").Value = -711.0001999999999
$ws.Range("N132").Value = -13051.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 419.5
$ws.Range("I26").Value = 419.5
$ws.Range("K26").Value = 419.5
$ws.Range("M26").Value = -89.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4702.048
$ws.Range("I32").Value = 4084.1052
$ws.Range("K32").Value = 4084.1052
$ws.Range("M32").Value = -3797.1052

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1885.5834
$ws.Range("I61").Value = 1885.5834
$ws.Range("K61").Value = 1885.5834
$ws.Range("M61").Value = -1673.5834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 45775
$ws.Range("J92").Value = 45775
$ws.Range("L92").Value = 45775
$ws.Range("N92").Value = -50767

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1885.5834
$ws.Range("I136").Value = 1885.5834
$ws.Range("K136").Value = 5656.7502
$ws.Range("M136").Value = -3106.7502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2577.4167
$ws.Range("I94").Value = 2189.875
$ws.Range("J94").Value = 3352.5
$ws.Range("K94").Value = 2189.875
$ws.Range("L94").Value = 3352.5
$ws.Range("M94").Value = -1738.875
$ws.Range("N94").Value = -4254.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2877.7222
$ws.Range("I99").Value = 1562.375
$ws.Range("J99").Value = 3930
$ws.Range("K99").Value = 1562.375
$ws.Range("L99").Value = 3930
$ws.Range("M99").Value = -64.375
$ws.Range("N99").Value = -6926

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3671.8
$ws.Range("I105").Value = 3428.1538
$ws.Range("K105").Value = 3428.1538
$ws.Range("M105").Value = -1681.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1542.5
$ws.Range("I134").Value = 1204.2
$ws.Range("K134").Value = 3612.6
$ws.Range("M134").Value = -1077.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51748

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -158736

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1791
$ws.Range("I105").Value = 1736.1538
$ws.Range("K105").Value = 1736.1538
$ws.Range("M105").Value = 10.84619999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1529.8889
$ws.Range("J122").Value = 1998
$ws.Range("L122").Value = 5994
$ws.Range("N122").Value = -10894

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 4.857143
$ws.Range("I12").Value = 3.8333333
$ws.Range("J12").Value = 11
$ws.Range("K12").Value = 11.4999999
$ws.Range("L12").Value = 33
$ws.Range("M12").Value = 161.5000001
$ws.Range("N12").Value = -379

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 150.6923
$ws.Range("J33").Value = 213.5
$ws.Range("L33").Value = 1281
$ws.Range("N33").Value = -1847

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 7670
$ws.Range("I94").Value = 1666.6666
$ws.Range("J94").Value = 10242.857
$ws.Range("K94").Value = 4999.9998
$ws.Range("L94").Value = 30728.571
$ws.Range("M94").Value = -4323.9998
$ws.Range("N94").Value = -32080.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1193.875
$ws.Range("J107").Value = 429.05264
$ws.Range("L107").Value = 1287.15792
$ws.Range("N107").Value = -5127.15792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3556.0833
$ws.Range("J129").Value = 5828.5713
$ws.Range("L129").Value = 17485.7139
$ws.Range("N129").Value = -27485.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4724.4165
$ws.Range("I137").Value = 1653.8334
$ws.Range("J137").Value = 7795
$ws.Range("K137").Value = 4961.5002
$ws.Range("L137").Value = 23385
$ws.Range("M137").Value = 138.4997999999996
$ws.Range("N137").Value = -33585

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 450.54166
$ws.Range("I97").Value = 427.8421
$ws.Range("J97").Value = 536.8
$ws.Range("K97").Value = 427.8421
$ws.Range("L97").Value = 536.8
$ws.Range("M97").Value = 68.15789999999998
$ws.Range("N97").Value = -1528.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5637.737
$ws.Range("I40").Value = 3686.7
$ws.Range("J40").Value = 7805.5557
$ws.Range("K40").Value = 3686.7
$ws.Range("L40").Value = 7805.5557
$ws.Range("M40").Value = -3550.7
$ws.Range("N40").Value = -8077.5557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2836.6875
$ws.Range("J46").Value = 2726.6667
$ws.Range("L46").Value = 2726.6667
$ws.Range("N46").Value = -3102.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3272.077
$ws.Range("J81").Value = 3876.8572
$ws.Range("L81").Value = 7753.7144
$ws.Range("N81").Value = -9875.714400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3272.077
$ws.Range("J84").Value = 3876.8572
$ws.Range("L84").Value = 38768.572
$ws.Range("N84").Value = -49376.572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 42549.832
$ws.Range("J124").Value = 42549.832
$ws.Range("L124").Value = 42549.832
$ws.Range("N124").Value = -52369.832

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1977.4445
$ws.Range("I126").Value = 1662.125
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 4986.375
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -2516.375
$ws.Range("N126").Value = -18440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5586.137
$ws.Range("I132").Value = 6617.1
$ws.Range("J132").Value = 1837.1818
$ws.Range("K132").Value = 19851.3
$ws.Range("L132").Value = 5511.5454
$ws.Range("M132").Value = -17321.3
$ws.Range("N132").Value = -10571.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 55000
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 100000
$ws.Range("M141").Value = -4820
$ws.Range("N141").Value = -110360
